$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-25 hold placeholder / misplaced values in the source file;
# clear them out fully before rewriting so no stray cells remain.
$ws.Range("A10:C25").Clear()

# Row 10
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Abordar os princípios básicos da termodinâmica de forma que os estudantes e futuros engenheiros tenham um entendimento claro e sólido sobre estes princípios. Apresentar diversos exemplos de engenharia do mundo real e de como a termodinâmica é aplicada na prática de engenharia. Enfatizar a compreensão da termodinâmica baseada na Física e em argumentos físicos, buscando incentivar o entendimento mais profundo da termodinâmica.'
$ws.Range("C10").Value = 'Abordar os princípios básicos da termodinâmica de forma que os estudantes e futuros engenheiros tenham um entendimento claro e sólido sobre estes princípios. Apresentar diversos exemplos de engenharia do mundo real e de como a termodinâmica é aplicada na prática de engenharia. Enfatizar a compreensão da termodinâmica baseada na Física e em argumentos físicos, buscando incentivar o entendimento mais profundo da termodinâmica.'
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Range("A11").Value = 'Objectives:'
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Range("A12").Value = 'Docentes responsáveis:'

# Row 13
$ws.Range("B13").Value = '5840521 - Rosa Ana Conte'
$ws.Range("C13").Value = '5840521 - Rosa Ana Conte'

# Row 14
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = '1. Termodinâmica e Energia. 2. Importância das unidades e análise dimensional.3. Sistemas e volumes de controle. 4. Equipamentos domésticos e a Termodinâmica. 5. Propriedades de um sistema: estados termodinâmicos e equilíbrio. 6. Eficiência na conversão de energia. 7. Processos e ciclos térmicos. 8. Termodinâmica e o meio ambiente.'
$ws.Range("C14").Value = '1. Termodinâmica e Energia. 2. Importância das unidades e análise dimensional.3. Sistemas e volumes de controle. 4. Equipamentos domésticos e a Termodinâmica. 5. Propriedades de um sistema: estados termodinâmicos e equilíbrio. 6. Eficiência na conversão de energia. 7. Processos e ciclos térmicos. 8. Termodinâmica e o meio ambiente.'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '1. Termodinâmica e Energia: formas de energia e transferência de energia por calor e trabalho; formas mecânicas de trabalho. 2. Sistema de Unidades e Análise Dimensional: importância na engenharia de máquinas. 3. Sistemas e volumes de controle: dispositivos ativos e passivos. 4. Propriedades de um sistema. Estados e equilíbrio: diagramas de propriedades para processos com mudança de fase; equilíbrio de estado do gás ideal; fator de compressibilidade; pressão de vapor e pressão de equilíbrio; calores específicos. 5. Balanço de energia em sistemas fechados e em volumes de controle: trabalho de fluxo e energia de escoamento de um fluido; regime permanente e transiente. 6. Máquinas térmicas e refrigeradores e a 2ª. Lei da Termodinâmica: princípios e ciclos de Carnot; entropia e variação de entropia em sólidos, líquidos e gases. 7. Eficiência na conversão de energia. Eficiência térmica. Eficiência de máquinas. Eficiência isoentrópica em dispositivos com escoamento em regime permanente. Balanço de entropia. 8. Processo e ciclos: Ciclos de potência a gás: Otto, Diesel, Stirling, Ericsson, Brayton e suas variações. Ciclos de potência a vapor e ciclos combinados gás-vapor: Rankine ideal; afastamento da condição ideal; eficiência do ciclo Rankine com e sem modificações; cogeração. Ciclos de refrigeração e sistemas de bombas de calor: sistemas a gás e por absorção. 9. Economia de energia: benefícios ao meio ambiente.'
$ws.Range("C16").Value = '1. Termodinâmica e Energia: formas de energia e transferência de energia por calor e trabalho; formas mecânicas de trabalho. 2. Sistema de Unidades e Análise Dimensional: importância na engenharia de máquinas. 3. Sistemas e volumes de controle: dispositivos ativos e passivos. 4. Propriedades de um sistema. Estados e equilíbrio: diagramas de propriedades para processos com mudança de fase; equilíbrio de estado do gás ideal; fator de compressibilidade; pressão de vapor e pressão de equilíbrio; calores específicos. 5. Balanço de energia em sistemas fechados e em volumes de controle: trabalho de fluxo e energia de escoamento de um fluido; regime permanente e transiente. 6. Máquinas térmicas e refrigeradores e a 2ª. Lei da Termodinâmica: princípios e ciclos de Carnot; entropia e variação de entropia em sólidos, líquidos e gases. 7. Eficiência na conversão de energia. Eficiência térmica. Eficiência de máquinas. Eficiência isoentrópica em dispositivos com escoamento em regime permanente. Balanço de entropia. 8. Processo e ciclos: Ciclos de potência a gás: Otto, Diesel, Stirling, Ericsson, Brayton e suas variações. Ciclos de potência a vapor e ciclos combinados gás-vapor: Rankine ideal; afastamento da condição ideal; eficiência do ciclo Rankine com e sem modificações; cogeração. Ciclos de refrigeração e sistemas de bombas de calor: sistemas a gás e por absorção. 9. Economia de energia: benefícios ao meio ambiente.'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = 'Syllabus:'
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A18").Value = 'Avaliação:'

# Row 19
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'Serão realizadas 2 avaliações, com questões abrangendo problemas práticos e conceituais. A 1a. avaliação terá peso 1 e a 2a. avaliação terá peso 2. A nota será a média ponderada das 2 avaliações.'
$ws.Range("C19").Value = 'Serão realizadas 2 avaliações, com questões abrangendo problemas práticos e conceituais. A 1a. avaliação terá peso 1 e a 2a. avaliação terá peso 2. A nota será a média ponderada das 2 avaliações.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'Serão aplicadas duas avaliações escritas (P1, com peso 1 e P2, com peso 2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (P1 + P2)/3.'
$ws.Range("C20").Value = 'Serão aplicadas duas avaliações escritas (P1, com peso 1 e P2, com peso 2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (P1 + P2)/3.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Range("C21").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = '1.Çengel, Y.A.; Boles, M.A. Thermodynamics An Engineering Approach, 6th ed., New York: McGraw Hill, 20082.Borgnakke, C; Sonntag, R.E. Fundamentos da termodinâmica, São Paulo: Blucher, 20133.Moran, M. J., Shapiro, H. N., Munson, B. R. & DeWitt, D. P. – Introdução à Engenharia de Sistemas Térmicos – LTC.4.Potter, M. C. & Scott, E. P. – Ciências Térmicas: Termodinâmica, Mecânica dos Fluidos e Transmissão de Calor – Thomson.5.Moran, M. J., Shapiro, H. N., Boettner, D. D. & Bailey, M. B. – Princípios de Termodinâmica para Engenharia – 7ª ed., LTC.6.Potter, M. C. & Scott, E. P. – Termodinâmica – Thomson.7.J.H. Keenan. Gas Tables: Thermodynamics Properties of Air Products of Combustion and Component Gases Compressible Flow Functions. John Wiley, 1980'
$ws.Range("C22").Value = '1.Çengel, Y.A.; Boles, M.A. Thermodynamics An Engineering Approach, 6th ed., New York: McGraw Hill, 20082.Borgnakke, C; Sonntag, R.E. Fundamentos da termodinâmica, São Paulo: Blucher, 20133.Moran, M. J., Shapiro, H. N., Munson, B. R. & DeWitt, D. P. – Introdução à Engenharia de Sistemas Térmicos – LTC.4.Potter, M. C. & Scott, E. P. – Ciências Térmicas: Termodinâmica, Mecânica dos Fluidos e Transmissão de Calor – Thomson.5.Moran, M. J., Shapiro, H. N., Boettner, D. D. & Bailey, M. B. – Princípios de Termodinâmica para Engenharia – 7ª ed., LTC.6.Potter, M. C. & Scott, E. P. – Termodinâmica – Thomson.7.J.H. Keenan. Gas Tables: Thermodynamics Properties of Air Products of Combustion and Component Gases Compressible Flow Functions. John Wiley, 1980'
$ws.Rows.Item(22).RowHeight = 120

# Row 23
$ws.Range("A23").Value = 'Requisitos:'

# Row 24
$ws.Range("B24").Value = 'LOB1004 -  Cálculo II  (Requisito)
'
$ws.Range("C24").Value = 'LOB1004 -  Cálculo II  (Requisito)
'
$ws.Rows.Item(24).RowHeight = 30

# Row 25
$ws.Range("B25").Value = 'LOB1019 -  Física II  (Requisito)
'
$ws.Range("C25").Value = 'LOB1019 -  Física II  (Requisito)
'
$ws.Rows.Item(25).RowHeight = 30

# The A-column width definition originally spanned cols A:B (min=1,max=2);
# split it down to just column A to match the authored column layout.
$ws.Columns.Item(1).Hidden = $false

Write-Host "edit complete"
